# "range & mov effect x 10."
# Multiply the Range (Q) and Mov (R) columns of the monster stat tables
# by 10 (they were all populated with 1, and move to 10 - except the very
# first data row's Mov value, which becomes 30).

$wb = $excel.ActiveWorkbook

# --- Sheet "标准卡" (table1, data rows 4-311) ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("Q4:Q311").Value = 10
$ws1.Range("R5:R311").Value = 10
$ws1.Range("R4").Value = 30

$ws1.Activate()
$ws1.Range("R5").Select()

# --- Sheet "隐藏卡" (table2, data rows 4-9) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("Q4:Q9").Value = 10
$ws2.Range("R4:R9").Value = 10

$ws2.Activate()
$ws2.Range("R6").Select()
